# Re-sorts the comma-separated list of names/emails in column G
# ("Recorded By") of the active worksheet into case-sensitive
# (ASCII) ascending order, e.g.
#   "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
#   "system, System, backup@backdoor.com"  -> "System, backup@backdoor.com, system"

# Helper: sort a list of strings using true ordinal (case-sensitive, ASCII)
# ordering. Note: this runtime's Sort-Object / -ceq / -clt operators behave
# case-INsensitively, so we must compare using the .NET string .CompareTo()
# method (which IS ordinal/case-sensitive) and sort manually.
function Sort-Ordinal($list) {
    $arr = @($list)
    $n = $arr.Count
    for ($i = 0; $i -lt $n; $i++) {
        for ($j = 0; $j -lt ($n - $i - 1); $j++) {
            if ($arr[$j].CompareTo($arr[$j + 1]) -gt 0) {
                $tmp = $arr[$j]
                $arr[$j] = $arr[$j + 1]
                $arr[$j + 1] = $tmp
            }
        }
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the used range so we know how many rows to touch.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Column G is the 7th column ("Recorded By").
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $value = $cell.Value()

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $text = $value.ToString()

    # Split on comma, trim whitespace from each part.
    $parts = $text.Split(",") | ForEach-Object { $_.Trim() }

    # Case-sensitive (ordinal/ASCII) sort, matching the original authoring tool's
    # behavior where uppercase letters sort before lowercase ones.
    $sorted = Sort-Ordinal $parts

    $newValue = [string]::Join(", ", $sorted)

    if ($newValue -ne $text) {
        $cell.Value = $newValue
    }
}
